$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 03:18:44"
$ws.Range("E3").Value = "2026-02-16 03:18:46"
$ws.Range("N3").Value = "-1.5 °C 2:54 TU"
$ws.Range("E4").Value = "2026-02-16 03:18:48"
$ws.Range("H4").Value = "'69%"
$ws.Range("E5").Value = "2026-02-16 03:18:51"
$ws.Range("G5").Value = "134 cm"
$ws.Range("I5").Value = "1.5 mm"
$ws.Range("E6").Value = "2026-02-16 03:18:53"
$ws.Range("H6").Value = "'84%"
$ws.Range("M6").Value = "8.1 °C 2:37 TU"
$ws.Range("N6").Value = "6.3 °C 2:57 TU"
$ws.Range("O6").Value = "7.0 °C"
$ws.Range("E7").Value = "2026-02-16 03:18:56"
$ws.Range("J7").Value = "1015.2 hPa"
$ws.Range("E8").Value = "2026-02-16 03:18:58"
$ws.Range("J8").Value = "1015.1 hPa"
$ws.Range("L8").Value = "59.8 km/h - 329º 2:44 TU"
$ws.Range("O8").Value = "9.2 °C"
$ws.Range("E9").Value = "2026-02-16 03:19:01"
$ws.Range("H9").Value = "'94%"
$ws.Range("N9").Value = "5.2 °C 2:59 TU"
$ws.Range("O9").Value = "5.6 °C"
$ws.Range("E10").Value = "2026-02-16 03:19:03"
$ws.Range("E11").Value = "2026-02-16 03:19:05"
$ws.Range("E12").Value = "2026-02-16 03:19:07"
$ws.Range("E13").Value = "2026-02-16 03:19:10"
$ws.Range("H13").Value = "'86%"
$ws.Range("J13").Value = "1018.5 hPa"
$ws.Range("E14").Value = "2026-02-16 03:19:11"
$ws.Range("M14").Value = "13.4 °C 2:32 TU"
$ws.Range("O14").Value = "12.7 °C"
$ws.Range("E15").Value = "2026-02-16 03:19:12"
$ws.Range("O15").Value = "5.8 °C"
$ws.Range("E16").Value = "2026-02-16 03:19:13"
$ws.Range("H16").Value = "'80%"
$ws.Range("M16").Value = "-0.5 °C 2:52 TU"
$ws.Range("O16").Value = "-1.0 °C"
$ws.Range("E17").Value = "2026-02-16 03:19:14"
$ws.Range("L17").Value = "37.4 km/h - 274º 2:40 TU"
$ws.Range("E18").Value = "2026-02-16 03:19:15"
$ws.Range("O18").Value = "4.5 °C"
$ws.Range("E19").Value = "2026-02-16 03:19:16"
$ws.Range("L19").Value = "9.4 km/h - 133º 2:53 TU"
$ws.Range("O19").Value = "3.4 °C"
$ws.Range("E20").Value = "2026-02-16 03:19:17"
$ws.Range("H20").Value = "'88%"
$ws.Range("N20").Value = "-1.5 °C 2:52 TU"
$ws.Range("E21").Value = "2026-02-16 03:19:18"
$ws.Range("H21").Value = "'79%"
$ws.Range("J21").Value = "1017.1 hPa"
$ws.Range("N21").Value = "3.7 °C 2:48 TU"
$ws.Range("O21").Value = "5.1 °C"
$ws.Range("E22").Value = "2026-02-16 03:19:19"
$ws.Range("N22").Value = "-6.5 °C 2:38 TU"
$ws.Range("E23").Value = "2026-02-16 03:19:21"
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = ""
$ws.Range("I23").Value = ""
$ws.Range("K23").Value = ""
$ws.Range("L23").Value = ""
$ws.Range("M23").Value = ""
$ws.Range("N23").Value = ""
$ws.Range("O23").Value = ""
$ws.Range("E24").Value = "2026-02-16 03:19:37"
$ws.Range("E25").Value = "2026-02-16 03:19:39"
$ws.Range("H25").Value = "'70%"
$ws.Range("E26").Value = "2026-02-16 03:19:42"
$ws.Range("E27").Value = "2026-02-16 03:19:44"
$ws.Range("E28").Value = "2026-02-16 03:19:46"
$ws.Range("H28").Value = "'90%"
$ws.Range("J28").Value = "1016.1 hPa"
$ws.Range("E29").Value = "2026-02-16 03:19:55"
$ws.Range("E30").Value = "2026-02-16 03:19:59"
$ws.Range("J30").Value = "1014.9 hPa"
$ws.Range("M30").Value = "7.7 °C 2:48 TU"
$ws.Range("O30").Value = "6.8 °C"
$ws.Range("E31").Value = "2026-02-16 03:20:01"
$ws.Range("J31").Value = "1013.4 hPa"
$ws.Range("L31").Value = "75.6 km/h - 328º 2:42 TU"
$ws.Range("O31").Value = "14.1 °C"
$ws.Range("E32").Value = "2026-02-16 03:20:04"
$ws.Range("E33").Value = "2026-02-16 03:20:06"
$ws.Range("H33").Value = "'69%"
$ws.Range("N33").Value = "4.2 °C 2:59 TU"
$ws.Range("O33").Value = "5.4 °C"
$ws.Range("E34").Value = "2026-02-16 03:20:09"
$ws.Range("N34").Value = "2.4 °C 2:34 TU"
$ws.Range("O34").Value = "3.3 °C"
$ws.Range("E35").Value = "2026-02-16 03:20:11"
$ws.Range("J35").Value = "1019.5 hPa"
$ws.Range("M35").Value = "7.0 °C 2:36 TU"
$ws.Range("E36").Value = "2026-02-16 03:20:14"
$ws.Range("H36").Value = "'87%"
$ws.Range("N36").Value = "5.4 °C 2:49 TU"
$ws.Range("O36").Value = "7.2 °C"
$ws.Range("E37").Value = "2026-02-16 03:20:16"
$ws.Range("H37").Value = "'93%"
$ws.Range("N37").Value = "1.8 °C 2:59 TU"
$ws.Range("E38").Value = "2026-02-16 03:20:19"
$ws.Range("H38").Value = "'92%"
$ws.Range("N38").Value = "4.7 °C 2:45 TU"
$ws.Range("O38").Value = "5.7 °C"
$ws.Range("E39").Value = "2026-02-16 03:20:21"
$ws.Range("M39").Value = "0.3 °C 2:59 TU"
$ws.Range("E40").Value = "2026-02-16 03:20:24"
$ws.Range("H40").Value = "'93%"
$ws.Range("N40").Value = "2.3 °C 2:59 TU"
$ws.Range("O40").Value = "3.4 °C"
$ws.Range("E41").Value = "2026-02-16 03:20:26"
$ws.Range("H41").Value = "'50%"
$ws.Range("J41").Value = "1016.2 hPa"
$ws.Range("E42").Value = "2026-02-16 03:20:29"
$ws.Range("O42").Value = "6.5 °C"
$ws.Range("E43").Value = "2026-02-16 03:20:31"
$ws.Range("N43").Value = "2.4 °C 2:37 TU"
$ws.Range("O43").Value = "3.7 °C"
$ws.Range("E44").Value = "2026-02-16 03:20:34"
$ws.Range("E45").Value = "2026-02-16 03:20:37"
$ws.Range("J45").Value = "1020.1 hPa"
$ws.Range("E46").Value = "2026-02-16 03:20:39"
$ws.Range("M46").Value = "13.2 °C 2:52 TU"
$ws.Range("O46").Value = "12.4 °C"
